# Applies the "optimized model for speed and size" edit:
# Adds a new worksheet "R. benchmark" with a small benchmark results table,
# and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "R. benchmark"

# Fill column A first (top to bottom), matching the order the strings were
# originally authored in, then the remaining columns/rows.
$newSheet.Range("A1").Value = "unigram"
$newSheet.Range("A2").Value = "top 500"
$newSheet.Range("A3").Value = "top 500"

$newSheet.Range("B1").Value = "alpha"
$newSheet.Range("B2").Value = 0.15
$newSheet.Range("B3").Value = 0.4

$newSheet.Range("C1").Value = "numer per root"
$newSheet.Range("C2").Value = "all"
$newSheet.Range("C3").Value = "all"

$resultText = "Overall top-3 score:     14.68 %`nOverall top-1 precision: 10.66 %`nOverall top-3 precision: 17.72 %`nAverage runtime:         910.16 msec`nNumber of predictions:   448`nTotal memory used:       2319.75 MB"
$newSheet.Range("D3").Value = $resultText
$newSheet.Range("D3").WrapText = $true
$newSheet.Rows.Item(3).RowHeight = 86.4

$newSheet.Range("D1").Value = "result"

# Column widths to roughly match the diff (bestFit width for C, custom for D)
$newSheet.Columns.Item("C").ColumnWidth = 13.44140625
$newSheet.Columns.Item("D").ColumnWidth = 49.44140625

# Selection matching diff (activeCell B3)
$newSheet.Range("B3").Select()

# Make new sheet the active (selected/visible) tab and deselect tab on Sheet2.
$newSheet.Activate()
